$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "35.226.62"
$cD.Style = $styleD
$ws.Range("E2").Value = "  +1.49%  "
$cD = $ws.Range("D3")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.885.21"
$cD.Style = $styleD
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.15%  "
$cD = $ws.Range("D5")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "244.86"
$cD.Style = $styleD
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  -0.05%  "
$cD = $ws.Range("D8")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "42.60"
$cD.Style = $styleD
$ws.Range("E8").Value = "  +3.78%  "
$cD = $ws.Range("D9")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.353"
$cD.Style = $styleD
$ws.Range("E9").Value = "  +3.98%  "
$cD = $ws.Range("D10")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "54.90"
$cD.Style = $styleD
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("E11").Value = "  +2.41%  "
$cD = $ws.Range("D12")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0983"
$cD.Style = $styleD
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +8.65%  "
$cD = $ws.Range("D14")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.775"
$cD.Style = $styleD
$ws.Range("E14").Value = "  +10.43%  "
$cD = $ws.Range("D15")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.158.18"
$cD.Style = $styleD
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("E16").Value = "  +3.51%  "
$cD = $ws.Range("D17")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.885.73"
$cD.Style = $styleD
$ws.Range("E17").Value = "  +1.11%  "
$cD = $ws.Range("D18")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "35.234.38"
$cD.Style = $styleD
$ws.Range("E18").Value = "  +1.53%  "
$cD = $ws.Range("D19")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "73.04"
$cD.Style = $styleD
$ws.Range("E19").Value = "  +1.95%  "
$cD = $ws.Range("D20")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0₃0821"
$cD.Style = $styleD
$ws.Range("E20").Value = "  +2.31%  "
$cD = $ws.Range("D21")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "243.35"
$cD.Style = $styleD
$ws.Range("E21").Value = "  +1.12%  "
$cD = $ws.Range("D22")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "12.74"
$cD.Style = $styleD
$ws.Range("E22").Value = "  +2.62%  "
$cD = $ws.Range("D23")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.14"
$cD.Style = $styleD
$ws.Range("E23").Value = "  +6.65%  "
$cD = $ws.Range("D24")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.65"
$cD.Style = $styleD
$ws.Range("E24").Value = "  +8.33%  "
$ws.Range("E25").Value = "  -0.03%  "
$cD = $ws.Range("D26")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "167.38"
$cD.Style = $styleD
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("E27").Value = "  -0.79%  "
$cD = $ws.Range("D28")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "8.45"
$cD.Style = $styleD
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("E30").Value = "  +1.52%  "
$cD = $ws.Range("D31")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "4.27"
$cD.Style = $styleD
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("E33").Value = "  +19.36%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -13.72%  "
$ws.Range("E37").Value = "  +4.05%  "
$cD = $ws.Range("D38")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.93"
$cD.Style = $styleD
$ws.Range("E38").Value = "  -0.03%  "
$cD = $ws.Range("D39")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0716"
$cD.Style = $styleD
$ws.Range("E39").Value = "  +7.22%  "
$ws.Range("E40").Value = "  +5.74%  "
$cD = $ws.Range("D41")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "97.70"
$cD.Style = $styleD
$ws.Range("E41").Value = "  +1.63%  "
$cD = $ws.Range("D42")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "17.06"
$cD.Style = $styleD
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("E43").Value = "  +1.35%  "
$cD = $ws.Range("D44")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.323.46"
$cD.Style = $styleD
$ws.Range("E44").Value = "  +4.15%  "
$cD = $ws.Range("D45")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "13.15"
$cD.Style = $styleD
$ws.Range("E45").Value = "  +11.57%  "
$cD = $ws.Range("D46")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.34"
$cD.Style = $styleD
$ws.Range("E46").Value = "  +2.82%  "
$cD = $ws.Range("D47")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.0811"
$cD.Style = $styleD
$ws.Range("E47").Value = "  +0.73%  "
$cD = $ws.Range("D48")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.41"
$cD.Style = $styleD
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +0.33%  "
$cD = $ws.Range("D51")
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.057.82"
$cD.Style = $styleD
$ws.Range("E51").Value = "  +0.77%  "
